$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0.0
$ws.Range("I16").Value = 0.0
$ws.Range("J16").Value = 0.0
$ws.Range("K16").Value = 0.0
$ws.Range("L16").Value = 0.0
$ws.Range("M16").Value = $null
$ws.Range("H80").Value = 1500.65
$ws.Range("I80").Value = 419.66666
$ws.Range("J80").Value = 1963.9286
$ws.Range("K80").Value = 1258.99998
$ws.Range("L80").Value = 5891.7858
$ws.Range("M80").Value = -260.9999800000001
$ws.Range("N80").Value = -7887.7858
$ws.Range("H83").Value = 1500.65
$ws.Range("I83").Value = 419.66666
$ws.Range("J83").Value = 1963.9286
$ws.Range("K83").Value = 3776.99994
$ws.Range("L83").Value = 17675.3574
$ws.Range("M83").Value = 1215.00006
$ws.Range("N83").Value = -27659.3574
$ws.Range("H104").Value = 146.8
$ws.Range("I104").Value = 133.5
$ws.Range("J104").Value = 200.0
$ws.Range("K104").Value = 400.5
$ws.Range("L104").Value = 600.0
$ws.Range("M104").Value = 1346.5
$ws.Range("N104").Value = -4094.0
$ws.Range("H127").Value = 2493.1428
$ws.Range("I127").Value = 3009.25
$ws.Range("J127").Value = 1805.0
$ws.Range("K127").Value = 9027.75
$ws.Range("L127").Value = 5415.0
$ws.Range("M127").Value = -4067.75
$ws.Range("N127").Value = -15335.0
$ws.Range("H129").Value = 2130.2942
$ws.Range("I129").Value = 371.0
$ws.Range("J129").Value = 3089.9092
$ws.Range("K129").Value = 1113.0
$ws.Range("L129").Value = 9269.7276
$ws.Range("M129").Value = 3887.0
$ws.Range("N129").Value = -19269.7276
$ws.Range("H132").Value = 2003.7273
$ws.Range("I132").Value = 2003.7273
$ws.Range("J132").Value = 0.0
$ws.Range("K132").Value = 6011.1819
$ws.Range("L132").Value = 0.0
$ws.Range("M132").Value = -3481.1819
$ws.Range("H138").Value = 4939.6
$ws.Range("I138").Value = 2441.4285
$ws.Range("J138").Value = 5911.1113
$ws.Range("K138").Value = 7324.2855
$ws.Range("L138").Value = 17733.3339
$ws.Range("M138").Value = -2184.2855
$ws.Range("H141").Value = 2085.8572
$ws.Range("I141").Value = 2085.8572
$ws.Range("J141").Value = 0.0
$ws.Range("K141").Value = 6257.571599999999
$ws.Range("L141").Value = 0.0
$ws.Range("M141").Value = -1077.571599999999
$ws.Range("N141").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3670185.5
$ws.Range("H88").Value = 2607.1667
$ws.Range("I88").Value = 2718.6667
$ws.Range("J88").Value = 2495.6667
$ws.Range("K88").Value = 2718.6667
$ws.Range("L88").Value = 2495.6667
$ws.Range("M88").Value = -2312.6667
$ws.Range("H91").Value = 2607.1667
$ws.Range("I91").Value = 2718.6667
$ws.Range("J91").Value = 2495.6667
$ws.Range("K91").Value = 2718.6667
$ws.Range("L91").Value = 2495.6667
$ws.Range("M91").Value = -1314.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 495.0
$ws.Range("I12").Value = 495.0
$ws.Range("J12").Value = 0.0
$ws.Range("K12").Value = 495.0
$ws.Range("L12").Value = 0.0
$ws.Range("M12").Value = -327.0
$ws.Range("N12").Value = $null
$ws.Range("H92").Value = 38874.5
$ws.Range("I92").Value = 0.0
$ws.Range("J92").Value = 38874.5
$ws.Range("K92").Value = 0.0
$ws.Range("L92").Value = 38874.5
$ws.Range("N92").Value = -43866.5
$ws.Range("H105").Value = 2866.6667
$ws.Range("I105").Value = 3000.0
$ws.Range("J105").Value = 2600.0
$ws.Range("K105").Value = 3000.0
$ws.Range("L105").Value = 2600.0
$ws.Range("M105").Value = -1253.0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 858.2
$ws.Range("I31").Value = 858.2
$ws.Range("J31").Value = 0.0
$ws.Range("K31").Value = 858.2
$ws.Range("L31").Value = 0.0
$ws.Range("M31").Value = -563.2
$ws.Range("N31").Value = $null
$ws.Range("H34").Value = 858.2
$ws.Range("I34").Value = 858.2
$ws.Range("J34").Value = 0.0
$ws.Range("K34").Value = 858.2
$ws.Range("L34").Value = 0.0
$ws.Range("M34").Value = -656.2
$ws.Range("N34").Value = $null
$ws.Range("H58").Value = 1777.2
$ws.Range("I58").Value = 1777.2
$ws.Range("J58").Value = 0.0
$ws.Range("K58").Value = 1777.2
$ws.Range("L58").Value = 0.0
$ws.Range("M58").Value = -1574.2
$ws.Range("N58").Value = $null
$ws.Range("H99").Value = 3300.0
$ws.Range("I99").Value = 2625.0
$ws.Range("J99").Value = 6000.0
$ws.Range("K99").Value = 2625.0
$ws.Range("L99").Value = 6000.0
$ws.Range("M99").Value = -1127.0
$ws.Range("H126").Value = 3300.0
$ws.Range("I126").Value = 2625.0
$ws.Range("J126").Value = 6000.0
$ws.Range("K126").Value = 7875.0
$ws.Range("L126").Value = 18000.0
$ws.Range("M126").Value = -5405.0
$ws.Range("H134").Value = 2157.2222
$ws.Range("I134").Value = 2249.2173
$ws.Range("J134").Value = 1628.25
$ws.Range("K134").Value = 6747.651899999999
$ws.Range("L134").Value = 4884.75
$ws.Range("M134").Value = -4212.651899999999
$ws.Range("N134").Value = -9954.75
$ws.Range("H136").Value = 1777.2
$ws.Range("I136").Value = 1777.2
$ws.Range("J136").Value = 0.0
$ws.Range("K136").Value = 5331.6
$ws.Range("L136").Value = 0.0
$ws.Range("M136").Value = -2781.6
$ws.Range("N136").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 15714482.0
$ws.Range("I4").Value = 18333546.0
$ws.Range("J4").Value = 100.0
$ws.Range("K4").Value = 55000638.0
$ws.Range("L4").Value = 300.0
$ws.Range("M4").Value = -55000526.0
$ws.Range("H127").Value = 1899.0
$ws.Range("I127").Value = 0.0
$ws.Range("J127").Value = 1899.0
$ws.Range("K127").Value = 0.0
$ws.Range("L127").Value = 5697.0
$ws.Range("N127").Value = -15617.0
$ws.Range("H129").Value = 1113515.2
$ws.Range("I129").Value = 2274.8
$ws.Range("J129").Value = 2502565.8
$ws.Range("K129").Value = 6824.400000000001
$ws.Range("L129").Value = 7507697.399999999
$ws.Range("M129").Value = -1824.400000000001
$ws.Range("N129").Value = -7517697.399999999
$ws.Range("H131").Value = 590695.8
$ws.Range("I131").Value = 1000.0
$ws.Range("J131").Value = 627551.8
$ws.Range("K131").Value = 3000.0
$ws.Range("L131").Value = 1882655.4
$ws.Range("M131").Value = 2040.0
$ws.Range("N131").Value = -1892735.4
$ws.Range("H140").Value = 9202.4
$ws.Range("I140").Value = 5515.0
$ws.Range("J140").Value = 10124.25
$ws.Range("K140").Value = 16545.0
$ws.Range("L140").Value = 30372.75
$ws.Range("M140").Value = -11365.0
$ws.Range("N140").Value = -40732.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 30030.0
$ws.Range("I52").Value = 30030.0
$ws.Range("J52").Value = 0.0
$ws.Range("K52").Value = 30030.0
$ws.Range("L52").Value = 0.0
$ws.Range("M52").Value = -29771.0
$ws.Range("H92").Value = 14220.0
$ws.Range("I92").Value = 0.0
$ws.Range("J92").Value = 14220.0
$ws.Range("K92").Value = 0.0
$ws.Range("L92").Value = 14220.0
$ws.Range("N92").Value = -17964.0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 894.0
$ws.Range("I22").Value = 894.0
$ws.Range("J22").Value = 0.0
$ws.Range("K22").Value = 894.0
$ws.Range("L22").Value = 0.0
$ws.Range("M22").Value = -599.0
$ws.Range("N22").Value = $null
$ws.Range("H27").Value = 894.0
$ws.Range("I27").Value = 894.0
$ws.Range("J27").Value = 0.0
$ws.Range("K27").Value = 894.0
$ws.Range("L27").Value = 0.0
$ws.Range("M27").Value = -787.0
$ws.Range("N27").Value = $null
$ws.Range("H46").Value = 2934.625
$ws.Range("I46").Value = 892.6667
$ws.Range("J46").Value = 4159.8
$ws.Range("K46").Value = 892.6667
$ws.Range("L46").Value = 4159.8
$ws.Range("M46").Value = -704.6667
$ws.Range("N46").Value = -4535.8
$ws.Range("H55").Value = 1109.8889
$ws.Range("I55").Value = 584.2857
$ws.Range("J55").Value = 2949.5
$ws.Range("K55").Value = 584.2857
$ws.Range("L55").Value = 2949.5
$ws.Range("M55").Value = -411.2857
$ws.Range("N55").Value = -3295.5
$ws.Range("H57").Value = 20041.0
$ws.Range("I57").Value = 20041.0
$ws.Range("J57").Value = 0.0
$ws.Range("K57").Value = 20041.0
$ws.Range("L57").Value = 0.0
$ws.Range("M57").Value = -19475.0
$ws.Range("N57").Value = $null
$ws.Range("H132").Value = 7371.8335
$ws.Range("I132").Value = 8595.889
$ws.Range("J132").Value = 3699.6667
$ws.Range("K132").Value = 25787.667
$ws.Range("L132").Value = 11099.0001
$ws.Range("M132").Value = -23257.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 19726.834
$ws.Range("I41").Value = 19715.5
$ws.Range("J41").Value = 19732.5
$ws.Range("K41").Value = 19715.5
$ws.Range("L41").Value = 19732.5
$ws.Range("M41").Value = -19325.5
